# Update cryptocurrency price/volume data per the commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.035.93"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.906.93"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7620"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.68"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3084"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.52"
$ws.Range("E9").Value = "  -6.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06903"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08012"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7556"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").Value = "1.907.48"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.262"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.77"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.193"
$ws.Range("E16").Value = "  +5.87%  "
$ws.Range("D17").Value = "30.039.39"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.05"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007746"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.43"
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "2.153.27"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.037"
$ws.Range("E24").Value = "  +5.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.318"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.44"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.86"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1316"
$ws.Range("E28").Value = "  +3.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.073"
$ws.Range("E29").Value = "  -2.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.344"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.049"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05457"
$ws.Range("E34").Value = "  +5.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.292"
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7376"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01945"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.767"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4455"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.97"
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.947"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8306"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.677"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.70"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.877"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "2.059.28"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.55"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1162"
